$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 18.375
$ws.Range("I6").Value = 18.375
$ws.Range("K6").Value = 55.125
$ws.Range("M6").Value = 56.875
$ws.Range("H12").Value = 1066.6666
$ws.Range("I12").Value = 725
$ws.Range("J12").Value = 1750
$ws.Range("K12").Value = 725
$ws.Range("L12").Value = 1750
$ws.Range("M12").Value = -555
$ws.Range("N12").Value = -2090
$ws.Range("H16").Value = 10822
$ws.Range("J16").Value = 13746.75
$ws.Range("L16").Value = 13746.75
$ws.Range("N16").Value = -14206.75
$ws.Range("H80").Value = 4466829.5
$ws.Range("I80").Value = 6252016.5
$ws.Range("J80").Value = 3862.875
$ws.Range("K80").Value = 18756049.5
$ws.Range("L80").Value = 11588.625
$ws.Range("M80").Value = -18755051.5
$ws.Range("N80").Value = -13584.625
$ws.Range("H83").Value = 4466829.5
$ws.Range("I83").Value = 6252016.5
$ws.Range("J83").Value = 3862.875
$ws.Range("K83").Value = 56268148.5
$ws.Range("L83").Value = 34765.875
$ws.Range("M83").Value = -56263156.5
$ws.Range("N83").Value = -44749.875
$ws.Range("H86").Value = 3737.348
$ws.Range("I86").Value = 2237
$ws.Range("J86").Value = 4393.75
$ws.Range("K86").Value = 2237
$ws.Range("L86").Value = 4393.75
$ws.Range("M86").Value = -1114
$ws.Range("N86").Value = -6639.75
$ws.Range("H88").Value = 5352.7896
$ws.Range("I88").Value = 3661.75
$ws.Range("J88").Value = 6582.636
$ws.Range("K88").Value = 3661.75
$ws.Range("L88").Value = 6582.636
$ws.Range("M88").Value = -3255.75
$ws.Range("N88").Value = -7394.636
$ws.Range("H89").Value = 3737.348
$ws.Range("I89").Value = 2237
$ws.Range("J89").Value = 4393.75
$ws.Range("K89").Value = 11185
$ws.Range("L89").Value = 21968.75
$ws.Range("M89").Value = -5569
$ws.Range("N89").Value = -33200.75
$ws.Range("H91").Value = 5352.7896
$ws.Range("I91").Value = 3661.75
$ws.Range("J91").Value = 6582.636
$ws.Range("K91").Value = 3661.75
$ws.Range("L91").Value = 6582.636
$ws.Range("M91").Value = -2257.75
$ws.Range("N91").Value = -9390.636
$ws.Range("H107").Value = 565.5
$ws.Range("I107").Value = 581.9
$ws.Range("K107").Value = 581.9
$ws.Range("M107").Value = 1338.1
$ws.Range("H112").Value = 1196.159
$ws.Range("J112").Value = 1208.7073
$ws.Range("L112").Value = 3626.1219
$ws.Range("N112").Value = -5842.1219
$ws.Range("H137").Value = 29413466
$ws.Range("J137").Value = 2709.4167
$ws.Range("L137").Value = 8128.250100000001
$ws.Range("N137").Value = -13228.2501
$ws.Range("H138").Value = 2766.257
$ws.Range("I138").Value = 1066.8462
$ws.Range("J138").Value = 3770.4546
$ws.Range("K138").Value = 3200.5386
$ws.Range("L138").Value = 11311.3638
$ws.Range("M138").Value = 1939.4614
$ws.Range("N138").Value = -21591.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1052.4237
$ws.Range("I32").Value = 1058.3392
$ws.Range("K32").Value = 1058.3392
$ws.Range("M32").Value = -771.3391999999999
$ws.Range("H74").Value = 2855.2222
$ws.Range("I74").Value = 1692.9333
$ws.Range("K74").Value = 1692.9333
$ws.Range("M74").Value = -818.9332999999999
$ws.Range("H77").Value = 2855.2222
$ws.Range("I77").Value = 1692.9333
$ws.Range("K77").Value = 8464.666499999999
$ws.Range("M77").Value = -4096.666499999999
$ws.Range("H97").Value = 914.75
$ws.Range("I97").Value = 914.75
$ws.Range("K97").Value = 914.75
$ws.Range("M97").Value = -418.75
$ws.Range("H122").Value = 1911.4166
$ws.Range("I122").Value = 1539.8182
$ws.Range("J122").Value = 5999
$ws.Range("K122").Value = 4619.4546
$ws.Range("L122").Value = 17997
$ws.Range("M122").Value = -2169.4546
$ws.Range("N122").Value = -22897

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2316.7778
$ws.Range("I86").Value = 2155
$ws.Range("J86").Value = 2883
$ws.Range("K86").Value = 2155
$ws.Range("L86").Value = 2883
$ws.Range("M86").Value = -1032
$ws.Range("N86").Value = -5129
$ws.Range("H89").Value = 2316.7778
$ws.Range("I89").Value = 2155
$ws.Range("J89").Value = 2883
$ws.Range("K89").Value = 10775
$ws.Range("L89").Value = 14415
$ws.Range("M89").Value = -5159
$ws.Range("N89").Value = -25647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3201.6667
$ws.Range("I31").Value = 1403.8
$ws.Range("K31").Value = 1403.8
$ws.Range("M31").Value = -1108.8
$ws.Range("H34").Value = 3201.6667
$ws.Range("I34").Value = 1403.8
$ws.Range("K34").Value = 1403.8
$ws.Range("M34").Value = -1201.8
$ws.Range("H107").Value = 821.4
$ws.Range("I107").Value = 689.5
$ws.Range("J107").Value = 909.3333
$ws.Range("K107").Value = 689.5
$ws.Range("L107").Value = 909.3333
$ws.Range("M107").Value = 1230.5
$ws.Range("N107").Value = -4749.3333
$ws.Range("H122").Value = 2172.125
$ws.Range("I122").Value = 1625.2858
$ws.Range("K122").Value = 4875.857400000001
$ws.Range("M122").Value = -2425.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 1742.5
$ws.Range("I100").Value = 1742.5
$ws.Range("K100").Value = 5227.5
$ws.Range("M100").Value = -4416.5
$ws.Range("H110").Value = 17521.215
$ws.Range("I110").Value = 4794.25
$ws.Range("J110").Value = 22612
$ws.Range("K110").Value = 14382.75
$ws.Range("L110").Value = 67836
$ws.Range("M110").Value = -10292.75
$ws.Range("N110").Value = -76016
$ws.Range("H117").Value = 5662.909
$ws.Range("J117").Value = 6532.4443
$ws.Range("L117").Value = 19597.3329
$ws.Range("N117").Value = -26481.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3530.2856
$ws.Range("I102").Value = 3368.6667
$ws.Range("J102").Value = 4500
$ws.Range("K102").Value = 3368.6667
$ws.Range("L102").Value = 4500
$ws.Range("M102").Value = -1746.6667
$ws.Range("N102").Value = -7744

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4834.8887
$ws.Range("I40").Value = 4834.8887
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4834.8887
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4698.8887
$ws.Range("N40").ClearContents()
$ws.Range("H55").Value = 314.42856
$ws.Range("I55").Value = 262.75
$ws.Range("J55").Value = 383.33334
$ws.Range("K55").Value = 262.75
$ws.Range("L55").Value = 383.33334
$ws.Range("M55").Value = -89.75
$ws.Range("N55").Value = -729.33334
$ws.Range("H122").Value = 4128.1665
$ws.Range("I122").Value = 3620.4666
$ws.Range("K122").Value = 10861.3998
$ws.Range("M122").Value = -8411.399800000001
$ws.Range("H132").Value = 3056.8572
$ws.Range("I132").Value = 2703.25
$ws.Range("J132").Value = 4188.4
$ws.Range("K132").Value = 8109.75
$ws.Range("L132").Value = 12565.2
$ws.Range("M132").Value = -5579.75
$ws.Range("N132").Value = -17625.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5669.8184
$ws.Range("I81").Value = 5669.8184
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 11339.6368
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -10278.6368
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 5669.8184
$ws.Range("I84").Value = 5669.8184
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 56698.184
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -51394.184
$ws.Range("N84").ClearContents()
